$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1219.4648
$ws.Range("I15").Value = 1219.4648
$ws.Range("K15").Value = 3658.3944
$ws.Range("M15").Value = -3489.3944
$ws.Range("H17").Value = 372070.12
$ws.Range("J17").Value = 372070.12
$ws.Range("L17").Value = 1116210.36
$ws.Range("N17").Value = -1116546.36
$ws.Range("H40").Value = 3409.75
$ws.Range("I40").Value = 2526.7144
$ws.Range("J40").Value = 3885.2307
$ws.Range("K40").Value = 2526.7144
$ws.Range("L40").Value = 3885.2307
$ws.Range("M40").Value = -2351.7144
$ws.Range("N40").Value = -4235.2307
$ws.Range("H92").Value = 616
$ws.Range("I92").Value = 577.6667
$ws.Range("K92").Value = 577.6667
$ws.Range("M92").Value = 670.3333
$ws.Range("H97").Value = 3505
$ws.Range("J97").Value = 3505
$ws.Range("L97").Value = 10515
$ws.Range("N97").Value = -11507
$ws.Range("H112").Value = 1507.5
$ws.Range("J112").Value = 1569.6428
$ws.Range("L112").Value = 4708.928400000001
$ws.Range("N112").Value = -6924.928400000001
$ws.Range("H127").Value = 806.2857
$ws.Range("I127").Value = 806.2857
$ws.Range("K127").Value = 2418.8571
$ws.Range("M127").Value = 2541.1429
$ws.Range("H132").Value = 2590.8572
$ws.Range("I132").Value = 1151.7142
$ws.Range("K132").Value = 3455.1426
$ws.Range("M132").Value = -925.1425999999997
$ws.Range("H137").Value = 3083.2778
$ws.Range("I137").Value = 1507.6923
$ws.Range("J137").Value = 7179.8
$ws.Range("K137").Value = 4523.0769
$ws.Range("L137").Value = 21539.4
$ws.Range("M137").Value = -1973.0769
$ws.Range("N137").Value = -26639.4
$ws.Range("H138").Value = 2110.4468
$ws.Range("I138").Value = 1464.6666
$ws.Range("J138").Value = 2632.0386
$ws.Range("K138").Value = 4393.9998
$ws.Range("L138").Value = 7896.1158
$ws.Range("M138").Value = 746.0002000000004
$ws.Range("N138").Value = -18176.1158

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3897
$ws.Range("I32").Value = 3581.3794
$ws.Range("J32").Value = 9999
$ws.Range("K32").Value = 3581.3794
$ws.Range("L32").Value = 9999
$ws.Range("M32").Value = -3294.3794
$ws.Range("N32").Value = -10573
$ws.Range("H61").Value = 2970.4707
$ws.Range("I61").Value = 3082.5
$ws.Range("J61").Value = 2909.3635
$ws.Range("K61").Value = 3082.5
$ws.Range("L61").Value = 2909.3635
$ws.Range("M61").Value = -2870.5
$ws.Range("N61").Value = -3333.3635
$ws.Range("H136").Value = 2970.4707
$ws.Range("I136").Value = 3082.5
$ws.Range("J136").Value = 2909.3635
$ws.Range("K136").Value = 9247.5
$ws.Range("L136").Value = 8728.0905
$ws.Range("M136").Value = -6697.5
$ws.Range("N136").Value = -13828.0905

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 12092.15
$ws.Range("J99").Value = 1011
$ws.Range("L99").Value = 1011
$ws.Range("N99").Value = -4007
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3406.6775
$ws.Range("I31").Value = 1403.95
$ws.Range("J31").Value = 7048
$ws.Range("K31").Value = 1403.95
$ws.Range("L31").Value = 7048
$ws.Range("M31").Value = -1108.95
$ws.Range("N31").Value = -7638
$ws.Range("H34").Value = 3406.6775
$ws.Range("I34").Value = 1403.95
$ws.Range("J34").Value = 7048
$ws.Range("K34").Value = 1403.95
$ws.Range("L34").Value = 7048
$ws.Range("M34").Value = -1201.95
$ws.Range("N34").Value = -7452
$ws.Range("H58").Value = 2143.1738
$ws.Range("I58").Value = 1818.1
$ws.Range("J58").Value = 2393.2307
$ws.Range("K58").Value = 1818.1
$ws.Range("L58").Value = 2393.2307
$ws.Range("M58").Value = -1615.1
$ws.Range("N58").Value = -2799.2307
$ws.Range("H107").Value = 2279.5806
$ws.Range("I107").Value = 1826.05
$ws.Range("J107").Value = 3104.182
$ws.Range("K107").Value = 1826.05
$ws.Range("L107").Value = 3104.182
$ws.Range("M107").Value = 93.95000000000005
$ws.Range("N107").Value = -6944.182
$ws.Range("H136").Value = 2143.1738
$ws.Range("I136").Value = 1818.1
$ws.Range("J136").Value = 2393.2307
$ws.Range("K136").Value = 5454.299999999999
$ws.Range("L136").Value = 7179.6921
$ws.Range("M136").Value = -2904.299999999999
$ws.Range("N136").Value = -12279.6921

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2794.3
$ws.Range("I75").Value = 698.6
$ws.Range("J75").Value = 4890
$ws.Range("K75").Value = 2095.8
$ws.Range("L75").Value = 14670
$ws.Range("M75").Value = -1097.8
$ws.Range("N75").Value = -16666
$ws.Range("H78").Value = 2794.3
$ws.Range("I78").Value = 698.6
$ws.Range("J78").Value = 4890
$ws.Range("K78").Value = 6287.400000000001
$ws.Range("L78").Value = 44010
$ws.Range("M78").Value = -1295.400000000001
$ws.Range("N78").Value = -53994
$ws.Range("H131").Value = 51971.65
$ws.Range("J131").Value = 2081.7693
$ws.Range("L131").Value = 6245.3079
$ws.Range("N131").Value = -16325.3079

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 37982.223
$ws.Range("I97").Value = 55423.418
$ws.Range("K97").Value = 55423.418
$ws.Range("M97").Value = -54927.418
$ws.Range("H126").Value = 2720.25
$ws.Range("I126").Value = 2720.25
$ws.Range("K126").Value = 8160.75
$ws.Range("M126").Value = -5690.75
$ws.Range("H132").Value = 4883.273
$ws.Range("I132").Value = 3506.1365
$ws.Range("K132").Value = 10518.4095
$ws.Range("M132").Value = -7988.4095

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 17574.96
$ws.Range("I100").Value = 3387.2856
$ws.Range("J100").Value = 34127.25
$ws.Range("K100").Value = 3387.2856
$ws.Range("L100").Value = 34127.25
$ws.Range("M100").Value = -2846.2856
$ws.Range("N100").Value = -35209.25
$ws.Range("H132").Value = 3784.818
$ws.Range("I132").Value = 3178.6572
$ws.Range("K132").Value = 9535.971600000001
$ws.Range("M132").Value = -7005.971600000001
$ws.Range("H136").Value = 3781
$ws.Range("I136").Value = 3053.25
$ws.Range("J136").Value = 4557.2666
$ws.Range("K136").Value = 9159.75
$ws.Range("L136").Value = 13671.7998
$ws.Range("M136").Value = -6609.75
$ws.Range("N136").Value = -18771.7998

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 30498
$ws.Range("J94").Value = 30498
$ws.Range("L94").Value = 30498
$ws.Range("N94").Value = -32300
$ws.Range("H100").Value = 5605.5713
$ws.Range("J100").Value = 2806.3333
$ws.Range("L100").Value = 5612.6666
$ws.Range("N100").Value = -6694.6666
$ws.Range("H132").Value = 5813546
$ws.Range("I132").Value = 5648135.5
$ws.Range("K132").Value = 16944406.5
$ws.Range("M132").Value = -16941876.5
